# Capital Typing research_data workbook update
# ----------------------------------------------
# The three "instruction" cells in column AP (rows 2-4) are rewritten so
# that they explicitly reference Capital Typing (www.capitaltyping.com)
# as the provider of the transcription services, instead of the older
# generic "Write 300-500 word article ..." / "Use the responses ..."
# wording. All other cell values are unaffected; any shared-string index
# churn in the saved XML is a natural side effect of rewriting these
# three strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - real estate lawyers instructions (AP2)
$ws.Range("AP2").Value = "Rewrite the content in 500 words. Refer to Capital Typing (www.capitaltyping.com) as a provider of real estate lawyers with legal Transcription services. "

# Row 3 - court reporters instructions (AP3)
$ws.Range("AP3").Value = "Rewrite the content in 500 words. Refer to Capital Typing (www.capitaltyping.com) as a provider of  legal Transcription services for court reporters."

# Row 4 - outsourcing transcription instructions (AP4)
$ws.Range("AP4").Value = "Use the responses provided to create a 500 word article for Capital Typing (www.capitaltyping.com) as a provider of  outsourcing transcription services."

# Reflect the author's final cell selection in the sheet view.
$ws.Range("AP7").Select()
